$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.759.12'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.807.74'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.95'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4305'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3671'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07195'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8598'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.76'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.923.15'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.574'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.329'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06882'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.72%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '80.23'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008879'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.17'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.794.14'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.188'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.50%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.152.76'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.61'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.856'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.21'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.214'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.89%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.899'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +15.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.05'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08921'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7510'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.161'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.404'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.769'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.52%  '
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.133'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +5.13%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.003'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05191'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01912'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5060'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1641'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.638'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.513'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +10.45%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.287'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.11%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '106.29'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.08%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.29'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.646'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.88%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06250'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4528'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.789'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.37%  '
